# Avances Controllers y DAOs.xlsx
# "Se agregan cambios a daos para plan de adquisiciones"
# Updates completion percentages on the Controllers and Daos sheets.

$wb = $excel.ActiveWorkbook

$wsControllers = $wb.Worksheets.Item("Controllers")
$wsDaos = $wb.Worksheets.Item("Daos")

# --- Sheet "Controllers" ---
$wsControllers.Range("C50").Value = 0.6

# --- Sheet "Daos" ---
$wsDaos.Range("C2").Value = 0.08
$wsDaos.Range("C11").Value = 1
$wsDaos.Range("C53").Value = 0.75
$wsDaos.Range("C96").Value = 0.1
$wsDaos.Range("C101").Value = 1

# Recalculate so the dependent COUNTIFS summary cells refresh their cached values
$excel.CalculateFullRebuild()

# --- Restore view/selection state (matches author's final navigation) ---
$wsDaos.Activate()
$wsDaos.Range("C54").Select()

$wsControllers.Activate()
$wsControllers.Range("C51").Select()
